# Rewrites the "Student Name / Marks / pass-fail" sheet into the new
# "Employe / Sales / Attendance / Bonus" sheet, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Carry formatting forward to the new/changed cells before we touch
# --- their contents, by copying format (not value) from cells that
# --- already wear the style we want.

# D2:D5 should end up with the old formula-column style (currently s=3,
# worn by C2). Grab it before C2's own style gets overwritten below.
$ws.Range("C2").Copy()
$ws.Range("D2:D5").PasteSpecial(-4122)

# C1 becomes a third header cell -> same style as A1/B1 (s=1).
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# C2:C5 becomes a third data column -> same style as A2:A5 / B2:B5 (s=2).
$ws.Range("B2:B5").Copy()
$ws.Range("C2:C5").PasteSpecial(-4122)

# --- Headers ---
$ws.Range("A1").Value = "Employe "
$ws.Range("B1").Value = "Sales "
$ws.Range("C1").Value = "Attendance "

# --- Data rows ---
$ws.Range("A2").Value = "Arjun"
$ws.Range("B2").Value = 55000
$ws.Range("C2").Value = 95

$ws.Range("A3").Value = "Simran "
$ws.Range("B3").Value = 65000
$ws.Range("C3").Value = 85

$ws.Range("A4").Value = "Rohit "
$ws.Range("B4").Value = 30000
$ws.Range("C4").Value = 90

$ws.Range("A5").Value = "Anjali "
$ws.Range("B5").Value = 25000
$ws.Range("C5").Value = 65

# --- Bonus formula column (shared across D2:D5) ---
$ws.Range("D2:D5").Formula = '=IF(AND(B2>40000,C2>=90%),"Bonus","No Bonus")'

# --- Drop the old 6th data row and the stray formatted cell at G9 ---
$ws.Range("A6:C6").Clear()
$ws.Range("G9").Clear()
